$p = $ppt.ActivePresentation

$s1 = $p.Slides.Item(1)
$t1 = $s1.Shapes.Item(1).TextFrame.TextRange
$t1.Text = "."
$t1.Text = "First slide"

$s3 = $p.Slides.Item(3)
$t3 = $s3.Shapes.Item(1).TextFrame.TextRange
$t3.Text = "."
$t3.Text = "Third slide"
